# Scheduled-runner market-data refresh for the Leve profit sheets.
# Updates the computed columns (H..N) on each class sheet with refreshed
# currentAveragePrice / LevePrice / LeveProfit figures. No formulas are
# involved - these are static snapshots written by the external runner -
# so we just overwrite the affected cell values (and clear the HQ-profit
# column where the HQ price has dropped to 0, since that column is omitted
# whenever there is no HQ price to compare against).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 195.3
$ws.Range("I4").Value = 195.3
$ws.Range("K4").Value = 195.3
$ws.Range("M4").Value = -81.30000000000001

$ws.Range("H62").Value = 1833
$ws.Range("I62").Value = 1666
$ws.Range("K62").Value = 1666
$ws.Range("M62").Value = -1042

$ws.Range("H65").Value = 1833
$ws.Range("I65").Value = 1666
$ws.Range("K65").Value = 8330
$ws.Range("M65").Value = -5210

$ws.Range("H86").Value = 78705620
$ws.Range("I86").Value = 225000800
$ws.Range("J86").Value = 5558027.5
$ws.Range("K86").Value = 225000800
$ws.Range("L86").Value = 5558027.5
$ws.Range("M86").Value = -224999677
$ws.Range("N86").Value = -5560273.5

$ws.Range("H89").Value = 78705620
$ws.Range("I89").Value = 225000800
$ws.Range("J89").Value = 5558027.5
$ws.Range("K89").Value = 1125004000
$ws.Range("L89").Value = 27790137.5
$ws.Range("M89").Value = -1124998384
$ws.Range("N89").Value = -27801369.5

$ws.Range("H111").Value = 12502138
$ws.Range("J111").Value = 3536.4
$ws.Range("L111").Value = 10609.2
$ws.Range("N111").Value = -16743.2

$ws.Range("H116").Value = 8626997
$ws.Range("I116").Value = 16669428
$ws.Range("J116").Value = 10106.214
$ws.Range("K116").Value = 16669428
$ws.Range("L116").Value = 10106.214
$ws.Range("M116").Value = -16665986
$ws.Range("N116").Value = -16990.214

$ws.Range("H132").Value = 779.2182
$ws.Range("I132").Value = 756.6111
$ws.Range("K132").Value = 2269.8333
$ws.Range("M132").Value = 260.1667000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 29413050
$ws.Range("I2").Value = 1034.3478
$ws.Range("K2").Value = 1034.3478
$ws.Range("M2").Value = -921.3478

$ws.Range("H44").Value = 67036.5
$ws.Range("J44").Value = 67036.5
$ws.Range("L44").Value = 67036.5
$ws.Range("N44").Value = -68012.5

$ws.Range("H45").Value = 4679.9
$ws.Range("I45").Value = 1503.7778
$ws.Range("J45").Value = 7278.5454
$ws.Range("K45").Value = 1503.7778
$ws.Range("L45").Value = 7278.5454
$ws.Range("M45").Value = -1126.7778
$ws.Range("N45").Value = -8032.5454

$ws.Range("H61").Value = 37042964
$ws.Range("I61").Value = 2152.353
$ws.Range("J61").Value = 100012344
$ws.Range("K61").Value = 2152.353
$ws.Range("L61").Value = 100012344
$ws.Range("M61").Value = -1940.353
$ws.Range("N61").Value = -100012768

$ws.Range("H106").Value = 53391.5
$ws.Range("J106").Value = 53391.5
$ws.Range("L106").Value = 53391.5
$ws.Range("N106").Value = -55915.5

$ws.Range("H116").Value = 29413050
$ws.Range("I116").Value = 1034.3478
$ws.Range("K116").Value = 1034.3478
$ws.Range("M116").Value = 1259.6522

$ws.Range("H130").Value = 19556.6
$ws.Range("J130").Value = 19556.6
$ws.Range("L130").Value = 19556.6
$ws.Range("N130").Value = -29596.6

$ws.Range("H136").Value = 37042964
$ws.Range("I136").Value = 2152.353
$ws.Range("J136").Value = 100012344
$ws.Range("K136").Value = 6457.059
$ws.Range("L136").Value = 300037032
$ws.Range("M136").Value = -3907.059
$ws.Range("N136").Value = -300042132

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 29413050
$ws.Range("I3").Value = 1034.3478
$ws.Range("K3").Value = 1034.3478
$ws.Range("M3").Value = -920.3478

$ws.Range("H86").Value = 7577851
$ws.Range("I86").Value = 10871334
$ws.Range("K86").Value = 10871334
$ws.Range("M86").Value = -10870211

$ws.Range("H89").Value = 7577851
$ws.Range("I89").Value = 10871334
$ws.Range("K89").Value = 54356670
$ws.Range("M89").Value = -54351054

$ws.Range("H134").Value = 6762355
$ws.Range("I134").Value = 13159455
$ws.Range("J134").Value = 9860.777
$ws.Range("K134").Value = 39478365
$ws.Range("L134").Value = 29582.331
$ws.Range("M134").Value = -39475830
$ws.Range("N134").Value = -34652.331

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -130
$ws.Range("N12").ClearContents()

$ws.Range("H13").Value = 4999.5
$ws.Range("J13").Value = 4999.5
$ws.Range("L13").Value = 4999.5
$ws.Range("N13").Value = -5277.5

$ws.Range("H28").Value = 41582
$ws.Range("J28").Value = 41582
$ws.Range("L28").Value = 41582
$ws.Range("N28").Value = -42072

$ws.Range("H35").Value = 932.3333
$ws.Range("I35").Value = 932.3333
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 932.3333
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -638.3333
$ws.Range("N35").ClearContents()

$ws.Range("H43").Value = 34557.668
$ws.Range("J43").Value = 34557.668
$ws.Range("L43").Value = 34557.668
$ws.Range("N43").Value = -34925.668

$ws.Range("H92").Value = 53642
$ws.Range("J92").Value = 53642
$ws.Range("L92").Value = 53642
$ws.Range("N92").Value = -58634

$ws.Range("H96").Value = 14972.182
$ws.Range("J96").Value = 14972.182
$ws.Range("L96").Value = 14972.182
$ws.Range("N96").Value = -20464.182

$ws.Range("H101").Value = 34557.668
$ws.Range("J101").Value = 34557.668
$ws.Range("L101").Value = 34557.668
$ws.Range("N101").Value = -41047.668

$ws.Range("H132").Value = 6862.5
$ws.Range("I132").Value = 4327.5884
$ws.Range("K132").Value = 12982.7652
$ws.Range("M132").Value = -10452.7652

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 12079.292
$ws.Range("I132").Value = 3539.5652
$ws.Range("J132").Value = 19935.84
$ws.Range("K132").Value = 31856.0868
$ws.Range("L132").Value = 179422.56
$ws.Range("M132").Value = -29326.0868
$ws.Range("N132").Value = -184482.56

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 38288
$ws.Range("J95").Value = 38288
$ws.Range("L95").Value = 38288
$ws.Range("N95").Value = -43780

$ws.Range("H100").Value = 49459.5
$ws.Range("J100").Value = 49459.5
$ws.Range("L100").Value = 49459.5
$ws.Range("N100").Value = -51623.5

$ws.Range("H122").Value = 9081431
$ws.Range("I122").Value = 9081431
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 27244293
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -27241843
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 5005.8823
$ws.Range("I132").Value = 2718.739
$ws.Range("J132").Value = 9788.091
$ws.Range("K132").Value = 8156.217000000001
$ws.Range("L132").Value = 29364.273
$ws.Range("M132").Value = -5626.217000000001
$ws.Range("N132").Value = -34424.273

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6116.5
$ws.Range("J22").Value = 6116.5
$ws.Range("L22").Value = 6116.5
$ws.Range("N22").Value = -6706.5

$ws.Range("H27").Value = 6116.5
$ws.Range("J27").Value = 6116.5
$ws.Range("L27").Value = 6116.5
$ws.Range("N27").Value = -6330.5

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H101").Value = 29381.8
$ws.Range("J101").Value = 29381.8
$ws.Range("L101").Value = 29381.8
$ws.Range("N101").Value = -35871.8

$ws.Range("H122").Value = 3225.4102
$ws.Range("I122").Value = 1995.5385
$ws.Range("J122").Value = 5685.154
$ws.Range("K122").Value = 5986.6155
$ws.Range("L122").Value = 17055.462
$ws.Range("M122").Value = -3536.6155
$ws.Range("N122").Value = -21955.462

$ws.Range("H132").Value = 10875998
$ws.Range("I132").Value = 22730040
$ws.Range("K132").Value = 68190120
$ws.Range("M132").Value = -68187590

$ws.Range("H136").Value = 10787.87
$ws.Range("I136").Value = 2924.9583
$ws.Range("K136").Value = 8774.874899999999
$ws.Range("M136").Value = -6224.874899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 214014.2
$ws.Range("I122").Value = 309959.62
$ws.Range("J122").Value = 6132.5
$ws.Range("K122").Value = 929878.86
$ws.Range("L122").Value = 18397.5
$ws.Range("M122").Value = -927428.86
$ws.Range("N122").Value = -23297.5

$ws.Range("H132").Value = 7303.143
$ws.Range("I132").Value = 7876.9443
$ws.Range("J132").Value = 6270.3
$ws.Range("K132").Value = 23630.8329
$ws.Range("L132").Value = 18810.9
$ws.Range("M132").Value = -21100.8329
$ws.Range("N132").Value = -23870.9

$ws.Range("H136").Value = 67341600
$ws.Range("I136").Value = 1000000000
$ws.Range("J136").Value = 723143.5600000001
$ws.Range("K136").Value = 3000000000
$ws.Range("L136").Value = 2169430.68
$ws.Range("M136").Value = -2999997450
$ws.Range("N136").Value = -2174530.68

$ws.Range("H138").Value = 68553.25
$ws.Range("J138").Value = 68553.25
$ws.Range("L138").Value = 68553.25
$ws.Range("N138").Value = -78833.25
